# ENW.xlsx — "added jiraids to ENW,ENWIAM,IAM xls"
#
# Row 2 (TCID "ENW001") gets a newly-linked Jira id (OPQA-1791) appended to
# its "Jira id" column, and a matching extra bullet appended to its
# "Description" column; the row grows from two wrapped lines to three, so
# the row height increases from 45 to 60. The previously-scrolled/selected
# cell (bottom of the sheet, C36) is reset back to C5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Append the new Jira id / description bullet to row 2 (TCID ENW001).
$ws.Range("B2").Value = "OPQA-1679||OPQA-3642||OPQA-1791"
$ws.Range("C2").Value = "Verify that the user is able to send the only one record at a time from article,Post,Patent view Pages||Verify that Non Market test group user who signed into the community application should be able to send a record to endnote.||Verify that User is able to sign-into EndNote Web with valid credentials"

# The extra bullet wraps onto another line, so the row is now taller.
$ws.Rows("2:2").RowHeight = 60

# Reset the saved view: selection back at C5, no special scroll anchor.
$ws.Range("C5").Select()
